# "Generate Report for Handoff"
# The localization status report moves from "In Translation" to
# "Ready for handoff": update the Status text and the associated
# timestamps on all three sheets, and widen the timestamp columns so the
# new values aren't clipped.

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Ready for handoff"
$newColWidth = 17.2159881591797

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew                  # zh-cn status
$wsOverview.Range("F2").Value = $statusNew                  # de-de status
$wsOverview.Range("G2").Value = "2016-09-01 18:46:15"       # Latest HO Xliff Generate Date
$wsOverview.Range("E1:F1").ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew                       # Status
$wsZhCn.Range("H2").Value = "2016-09-01 18:46:00"            # Latest Handoff Datetime
$wsZhCn.Range("C1").ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew                       # Status
$wsDeDe.Range("H2").Value = "2016-09-01 18:46:15"            # Latest Handoff Datetime
$wsDeDe.Range("C1").ColumnWidth = $newColWidth
